$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22224988
$ws.Range("I62").Value = 22224988
$ws.Range("K62").Value = 22224988
$ws.Range("M62").Value = -22224364
$ws.Range("H65").Value = 22224988
$ws.Range("I65").Value = 22224988
$ws.Range("K65").Value = 111124940
$ws.Range("M65").Value = -111121820
$ws.Range("H80").Value = 5433.05
$ws.Range("I80").Value = 379.16666
$ws.Range("J80").Value = 13013.875
$ws.Range("K80").Value = 1137.49998
$ws.Range("L80").Value = 39041.625
$ws.Range("M80").Value = -139.4999800000001
$ws.Range("N80").Value = -41037.625
$ws.Range("H83").Value = 5433.05
$ws.Range("I83").Value = 379.16666
$ws.Range("J83").Value = 13013.875
$ws.Range("K83").Value = 3412.49994
$ws.Range("L83").Value = 117124.875
$ws.Range("M83").Value = 1579.50006
$ws.Range("N83").Value = -127108.875
$ws.Range("H112").Value = 27212490
$ws.Range("J112").Value = 33615252
$ws.Range("L112").Value = 100845756
$ws.Range("N112").Value = -100847972
$ws.Range("H138").Value = 3944.3845
$ws.Range("I138").Value = 2054.2
$ws.Range("J138").Value = 4784.467
$ws.Range("K138").Value = 6162.599999999999
$ws.Range("L138").Value = 14353.401
$ws.Range("M138").Value = -1022.599999999999
$ws.Range("N138").Value = -24633.401

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17904.645
$ws.Range("I32").Value = 14432.018
$ws.Range("K32").Value = 14432.018
$ws.Range("M32").Value = -14145.018
$ws.Range("H45").Value = 141183.73
$ws.Range("J45").Value = 2450
$ws.Range("L45").Value = 2450
$ws.Range("N45").Value = -3204
$ws.Range("H61").Value = 196734
$ws.Range("I61").Value = 5296.4062
$ws.Range("J61").Value = 503034.16
$ws.Range("K61").Value = 5296.4062
$ws.Range("L61").Value = 503034.16
$ws.Range("M61").Value = -5084.4062
$ws.Range("N61").Value = -503458.16
$ws.Range("H110").Value = 9916.805
$ws.Range("I110").Value = 11228.333
$ws.Range("J110").Value = 4506.75
$ws.Range("K110").Value = 11228.333
$ws.Range("L110").Value = 4506.75
$ws.Range("M110").Value = -9183.333000000001
$ws.Range("N110").Value = -8596.75
$ws.Range("H136").Value = 196734
$ws.Range("I136").Value = 5296.4062
$ws.Range("J136").Value = 503034.16
$ws.Range("K136").Value = 15889.2186
$ws.Range("L136").Value = 1509102.48
$ws.Range("M136").Value = -13339.2186
$ws.Range("N136").Value = -1514202.48

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2352.8333
$ws.Range("I20").Value = 2602.6667
$ws.Range("J20").Value = 2103
$ws.Range("K20").Value = 2602.6667
$ws.Range("L20").Value = 2103
$ws.Range("M20").Value = -2355.6667
$ws.Range("N20").Value = -2597
$ws.Range("H105").Value = 2902.5
$ws.Range("I105").Value = 1601.9
$ws.Range("J105").Value = 9405.5
$ws.Range("K105").Value = 1601.9
$ws.Range("L105").Value = 9405.5
$ws.Range("M105").Value = 145.0999999999999
$ws.Range("N105").Value = -12899.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2343.875
$ws.Range("I99").Value = 1449.6666
$ws.Range("J99").Value = 2880.4
$ws.Range("K99").Value = 1449.6666
$ws.Range("L99").Value = 2880.4
$ws.Range("M99").Value = 48.33339999999998
$ws.Range("N99").Value = -5876.4
$ws.Range("H122").Value = 5567.25
$ws.Range("I122").Value = 5302.8887
$ws.Range("J122").Value = 5907.143
$ws.Range("K122").Value = 15908.6661
$ws.Range("L122").Value = 17721.429
$ws.Range("M122").Value = -13458.6661
$ws.Range("N122").Value = -22621.429
$ws.Range("H126").Value = 2343.875
$ws.Range("I126").Value = 1449.6666
$ws.Range("J126").Value = 2880.4
$ws.Range("K126").Value = 4348.9998
$ws.Range("L126").Value = 8641.200000000001
$ws.Range("M126").Value = -1878.9998
$ws.Range("N126").Value = -13581.2
$ws.Range("H134").Value = 8585964
$ws.Range("I134").Value = 9806847
$ws.Range("J134").Value = 1667625
$ws.Range("K134").Value = 29420541
$ws.Range("L134").Value = 5002875
$ws.Range("M134").Value = -29418006
$ws.Range("N134").Value = -5007945

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 778.625
$ws.Range("I92").Value = 280
$ws.Range("J92").Value = 849.8570999999999
$ws.Range("K92").Value = 840
$ws.Range("L92").Value = 2549.5713
$ws.Range("M92").Value = 408
$ws.Range("N92").Value = -5045.5713
$ws.Range("H113").Value = 192802.83
$ws.Range("I113").Value = 484.3889
$ws.Range("J113").Value = 294618.47
$ws.Range("K113").Value = 1453.1667
$ws.Range("L113").Value = 883855.4099999999
$ws.Range("M113").Value = 716.8333
$ws.Range("N113").Value = -888195.4099999999
$ws.Range("H131").Value = 2175004
$ws.Range("I131").Value = 9091679
$ws.Range("J131").Value = 1191.6571
$ws.Range("K131").Value = 27275037
$ws.Range("L131").Value = 3574.9713
$ws.Range("M131").Value = -27269997
$ws.Range("N131").Value = -13654.9713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5080.0938
$ws.Range("I70").Value = 4935.7144
$ws.Range("J70").Value = 5192.3887
$ws.Range("K70").Value = 4935.7144
$ws.Range("L70").Value = 5192.3887
$ws.Range("M70").Value = -4665.7144
$ws.Range("N70").Value = -5732.3887
$ws.Range("H73").Value = 5080.0938
$ws.Range("I73").Value = 4935.7144
$ws.Range("J73").Value = 5192.3887
$ws.Range("K73").Value = 4935.7144
$ws.Range("L73").Value = 5192.3887
$ws.Range("M73").Value = -3999.7144
$ws.Range("N73").Value = -7064.3887
$ws.Range("H102").Value = 3556
$ws.Range("I102").Value = 2993
$ws.Range("K102").Value = 2993
$ws.Range("M102").Value = -1371
$ws.Range("H103").Value = 17629.334
$ws.Range("J103").Value = 17629.334
$ws.Range("L103").Value = 17629.334
$ws.Range("N103").Value = -19973.334
$ws.Range("H113").Value = 66674264
$ws.Range("I113").Value = 125013180
$ws.Range("J113").Value = 1214.2858
$ws.Range("K113").Value = 125013180
$ws.Range("L113").Value = 1214.2858
$ws.Range("M113").Value = -125011010
$ws.Range("N113").Value = -5554.2858
$ws.Range("H132").Value = 5378346.5
$ws.Range("I132").Value = 6668149.5
$ws.Range("K132").Value = 20004448.5
$ws.Range("M132").Value = -20001918.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1458.6471
$ws.Range("I61").Value = 1253.2307
$ws.Range("K61").Value = 1253.2307
$ws.Range("M61").Value = -1051.2307
$ws.Range("H113").Value = 1458.6471
$ws.Range("I113").Value = 1253.2307
$ws.Range("K113").Value = 1253.2307
$ws.Range("M113").Value = 916.7692999999999
$ws.Range("H122").Value = 3396427.8
$ws.Range("I122").Value = 3972015
$ws.Range("J122").Value = 1669665.9
$ws.Range("K122").Value = 11916045
$ws.Range("L122").Value = 5008997.699999999
$ws.Range("M122").Value = -11913595
$ws.Range("N122").Value = -5013897.699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 446445.5
$ws.Range("I62").Value = 888888
$ws.Range("K62").Value = 888888
$ws.Range("M62").Value = -888264
$ws.Range("H65").Value = 446445.5
$ws.Range("I65").Value = 888888
$ws.Range("K65").Value = 4444440
$ws.Range("M65").Value = -4441320
$ws.Range("H122").Value = 1322.6428
$ws.Range("I122").Value = 1086.7142
$ws.Range("J122").Value = 1558.5714
$ws.Range("K122").Value = 3260.1426
$ws.Range("L122").Value = 4675.7142
$ws.Range("M122").Value = -810.1425999999997
$ws.Range("N122").Value = -9575.7142
